# ---------------------------------------------------------------------------
# edit.ps1 - Excel COM-interop script (Yojimbo_Profits sheets refresh)
#
# Updates the market-price-derived columns (H currentAveragePrice,
# I currentAveragePriceNQ, J currentAveragePriceHQ, K LevePriceNQ,
# L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ) on a handful of rows across
# all eight class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), mirroring a
# scheduled market-data refresh run.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 237134.31
$ws.Cells.Item(17, 10).Value = 237134.31
$ws.Cells.Item(17, 12).Value = 711402.9299999999
$ws.Cells.Item(17, 14).Value = -711738.9299999999

$ws.Cells.Item(87, 8).Value = 18275.455
$ws.Cells.Item(87, 10).Value = 18275.455
$ws.Cells.Item(87, 12).Value = 18275.455
$ws.Cells.Item(87, 14).Value = -20771.455

$ws.Cells.Item(90, 8).Value = 18275.455
$ws.Cells.Item(90, 10).Value = 18275.455
$ws.Cells.Item(90, 12).Value = 54826.36500000001
$ws.Cells.Item(90, 14).Value = -67306.36500000001

$ws.Cells.Item(129, 8).Value = 1005.6389
$ws.Cells.Item(129, 10).Value = 1028.9857
$ws.Cells.Item(129, 12).Value = 3086.9571
$ws.Cells.Item(129, 14).Value = -13086.9571

$ws.Cells.Item(132, 8).Value = 3131634.5
$ws.Cells.Item(132, 9).Value = 3683396.2
$ws.Cells.Item(132, 10).Value = 4984.3335
$ws.Cells.Item(132, 11).Value = 11050188.6
$ws.Cells.Item(132, 12).Value = 14953.0005
$ws.Cells.Item(132, 13).Value = -11047658.6
$ws.Cells.Item(132, 14).Value = -20013.0005

$ws.Cells.Item(137, 8).Value = 3217.0732
$ws.Cells.Item(137, 9).Value = 2871.9375
$ws.Cells.Item(137, 11).Value = 8615.8125
$ws.Cells.Item(137, 13).Value = -6065.8125

$ws.Cells.Item(138, 8).Value = 3047.3718
$ws.Cells.Item(138, 9).Value = 1270.6072
$ws.Cells.Item(138, 10).Value = 4042.36
$ws.Cells.Item(138, 11).Value = 3811.8216
$ws.Cells.Item(138, 12).Value = 12127.08
$ws.Cells.Item(138, 13).Value = 1328.1784
$ws.Cells.Item(138, 14).Value = -22407.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 28000
$ws.Cells.Item(123, 10).Value = 28000
$ws.Cells.Item(123, 12).Value = 28000
$ws.Cells.Item(123, 14).Value = -37800

$ws.Cells.Item(125, 8).Value = 40460
$ws.Cells.Item(125, 10).Value = 40460
$ws.Cells.Item(125, 12).Value = 40460
$ws.Cells.Item(125, 14).Value = -50300

$ws.Cells.Item(132, 8).Value = 45495.56
$ws.Cells.Item(132, 9).Value = 59752.35
$ws.Cells.Item(132, 10).Value = 15199.875
$ws.Cells.Item(132, 11).Value = 179257.05
$ws.Cells.Item(132, 12).Value = 45599.625
$ws.Cells.Item(132, 13).Value = -176727.05
$ws.Cells.Item(132, 14).Value = -50659.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 896
$ws.Cells.Item(107, 9).Value = 815.5454999999999
$ws.Cells.Item(107, 10).Value = 1092.6666
$ws.Cells.Item(107, 11).Value = 815.5454999999999
$ws.Cells.Item(107, 12).Value = 1092.6666
$ws.Cells.Item(107, 13).Value = 1104.4545
$ws.Cells.Item(107, 14).Value = -4932.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2553
$ws.Cells.Item(16, 9).Value = 3000
$ws.Cells.Item(16, 10).Value = 2404
$ws.Cells.Item(16, 11).Value = 3000
$ws.Cells.Item(16, 12).Value = 2404
$ws.Cells.Item(16, 13).Value = -2713
$ws.Cells.Item(16, 14).Value = -2978

$ws.Cells.Item(31, 8).Value = 31180.945
$ws.Cells.Item(31, 9).Value = 43014.883
$ws.Cells.Item(31, 10).Value = 3209.818
$ws.Cells.Item(31, 11).Value = 43014.883
$ws.Cells.Item(31, 12).Value = 3209.818
$ws.Cells.Item(31, 13).Value = -42719.883
$ws.Cells.Item(31, 14).Value = -3799.818

$ws.Cells.Item(34, 8).Value = 31180.945
$ws.Cells.Item(34, 9).Value = 43014.883
$ws.Cells.Item(34, 10).Value = 3209.818
$ws.Cells.Item(34, 11).Value = 43014.883
$ws.Cells.Item(34, 12).Value = 3209.818
$ws.Cells.Item(34, 13).Value = -42812.883
$ws.Cells.Item(34, 14).Value = -3613.818

$ws.Cells.Item(113, 8).Value = 2553
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 2404
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 2404
$ws.Cells.Item(113, 13).Value = -830
$ws.Cells.Item(113, 14).Value = -6744

$ws.Cells.Item(134, 8).Value = 7436.7354
$ws.Cells.Item(134, 9).Value = 6534.5
$ws.Cells.Item(134, 10).Value = 8451.75
$ws.Cells.Item(134, 11).Value = 19603.5
$ws.Cells.Item(134, 12).Value = 25355.25
$ws.Cells.Item(134, 13).Value = -17068.5
$ws.Cells.Item(134, 14).Value = -30425.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 4521.1904
$ws.Cells.Item(109, 9).Value = 1344.7
$ws.Cells.Item(109, 10).Value = 7408.909
$ws.Cells.Item(109, 11).Value = 4034.1
$ws.Cells.Item(109, 12).Value = 22226.727
$ws.Cells.Item(109, 13).Value = -2994.1
$ws.Cells.Item(109, 14).Value = -24306.727

$ws.Cells.Item(121, 8).Value = 41671292
$ws.Cells.Item(121, 9).Value = 83333580
$ws.Cells.Item(121, 10).Value = 9000
$ws.Cells.Item(121, 11).Value = 250000740
$ws.Cells.Item(121, 12).Value = 27000
$ws.Cells.Item(121, 13).Value = -249999430
$ws.Cells.Item(121, 14).Value = -29620

$ws.Cells.Item(134, 8).Value = 3802.889
$ws.Cells.Item(134, 9).Value = 1628.2354
$ws.Cells.Item(134, 10).Value = 7499.8
$ws.Cells.Item(134, 11).Value = 4884.706200000001
$ws.Cells.Item(134, 12).Value = 22499.4
$ws.Cells.Item(134, 13).Value = 185.2937999999995
$ws.Cells.Item(134, 14).Value = -32639.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2970.25
$ws.Cells.Item(80, 9).Value = 2816.6667
$ws.Cells.Item(80, 10).Value = 3062.4
$ws.Cells.Item(80, 11).Value = 2816.6667
$ws.Cells.Item(80, 12).Value = 3062.4
$ws.Cells.Item(80, 13).Value = -1818.6667
$ws.Cells.Item(80, 14).Value = -5058.4

$ws.Cells.Item(83, 8).Value = 2970.25
$ws.Cells.Item(83, 9).Value = 2816.6667
$ws.Cells.Item(83, 10).Value = 3062.4
$ws.Cells.Item(83, 11).Value = 14083.3335
$ws.Cells.Item(83, 12).Value = 15312
$ws.Cells.Item(83, 13).Value = -9091.333500000001
$ws.Cells.Item(83, 14).Value = -25296

$ws.Cells.Item(97, 8).Value = 769.9231
$ws.Cells.Item(97, 9).Value = 769.9231
$ws.Cells.Item(97, 11).Value = 769.9231
$ws.Cells.Item(97, 13).Value = -273.9231

$ws.Cells.Item(102, 8).Value = 1270.95
$ws.Cells.Item(102, 9).Value = 1251.0555
$ws.Cells.Item(102, 10).Value = 1450
$ws.Cells.Item(102, 11).Value = 1251.0555
$ws.Cells.Item(102, 12).Value = 1450
$ws.Cells.Item(102, 13).Value = 370.9445000000001
$ws.Cells.Item(102, 14).Value = -4694

$ws.Cells.Item(112, 8).Value = 9000
$ws.Cells.Item(112, 10).Value = 9000
$ws.Cells.Item(112, 12).Value = 9000
$ws.Cells.Item(112, 14).Value = -11216

$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 40572
$ws.Cells.Item(36, 10).Value = 40572
$ws.Cells.Item(36, 12).Value = 40572
$ws.Cells.Item(36, 14).Value = -41696

$ws.Cells.Item(88, 8).Value = 31500
$ws.Cells.Item(88, 10).Value = 31500
$ws.Cells.Item(88, 12).Value = 31500
$ws.Cells.Item(88, 14).Value = -32356

$ws.Cells.Item(91, 8).Value = 31500
$ws.Cells.Item(91, 10).Value = 31500
$ws.Cells.Item(91, 12).Value = 31500
$ws.Cells.Item(91, 14).Value = -34464

$ws.Cells.Item(125, 8).Value = 38000
$ws.Cells.Item(125, 10).Value = 38000
$ws.Cells.Item(125, 12).Value = 38000
$ws.Cells.Item(125, 14).Value = -47840

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 11500
$ws.Cells.Item(63, 10).Value = 11500
$ws.Cells.Item(63, 12).Value = 11500
$ws.Cells.Item(63, 14).Value = -12748

$ws.Cells.Item(66, 8).Value = 11500
$ws.Cells.Item(66, 10).Value = 11500
$ws.Cells.Item(66, 12).Value = 34500
$ws.Cells.Item(66, 14).Value = -40740

$ws.Cells.Item(69, 8).Value = 10000
$ws.Cells.Item(69, 10).Value = 10000
$ws.Cells.Item(69, 12).Value = 10000
$ws.Cells.Item(69, 14).Value = -11498

$ws.Cells.Item(72, 8).Value = 10000
$ws.Cells.Item(72, 10).Value = 10000
$ws.Cells.Item(72, 12).Value = 30000
$ws.Cells.Item(72, 14).Value = -37488

$ws.Cells.Item(123, 8).Value = 16916.666
$ws.Cells.Item(123, 10).Value = 16916.666
$ws.Cells.Item(123, 12).Value = 16916.666
$ws.Cells.Item(123, 14).Value = -26716.666

$ws.Cells.Item(125, 8).Value = 39759.312
$ws.Cells.Item(125, 10).Value = 39759.312
$ws.Cells.Item(125, 12).Value = 39759.312
$ws.Cells.Item(125, 14).Value = -49599.312
